# Apply updated crypto price/volume figures (and the Polygon/Polkadot
# and TrustWalletToken/VeChain row swaps) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe (escaped with a backtick inside the double-quoted
# string) forces Excel to store numeric-looking values such as "213.63"
# as text, matching the original inlineStr cell content exactly.

$ws.Range("D2").Value = "27.920.09"
$ws.Range("E2").Value = "  +1.49%  "

$ws.Range("D3").Value = "1.641.65"
$ws.Range("E3").Value = "  +1.20%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "`'213.63"
$ws.Range("E5").Value = "  +1.04%  "

$ws.Range("D6").Value = "`'0.524"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "`'23.88"
$ws.Range("E8").Value = "  +3.68%  "

$ws.Range("E9").Value = "  +0.78%  "

$ws.Range("E10").Value = "  +0.83%  "

$ws.Range("D11").Value = "`'0.0877"
$ws.Range("E11").Value = "  -0.14%  "

$ws.Range("D12").Value = "1.875.01"
$ws.Range("E12").Value = "  +1.19%  "

$ws.Range("D13").Value = "1.636.18"
$ws.Range("E13").Value = "  +0.85%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "`'4.08"
$ws.Range("E14").Value = "  +1.13%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "`'0.573"
$ws.Range("E15").Value = "  +4.57%  "

$ws.Range("E16").Value = "  +1.29%  "

$ws.Range("D17").Value = "27.918.13"
$ws.Range("E17").Value = "  +1.57%  "

$ws.Range("D18").Value = "`'231.13"
$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("E19").Value = "  +1.33%  "

$ws.Range("D20").Value = "`'7.61"
$ws.Range("E20").Value = "  +1.21%  "

$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").Value = "`'11.19"
$ws.Range("E22").Value = "  +7.96%  "

$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("D24").Value = "`'2.06"
$ws.Range("E24").Value = "  -0.72%  "

$ws.Range("D25").Value = "`'152.18"
$ws.Range("E25").Value = "  +2.14%  "

$ws.Range("E26").Value = "  +1.01%  "

$ws.Range("E27").Value = "  +0.78%  "

$ws.Range("D28").Value = "`'15.74"
$ws.Range("E28").Value = "  +1.40%  "

$ws.Range("E30").Value = "  +1.13%  "

$ws.Range("E31").Value = "  +0.75%  "

$ws.Range("E32").Value = "  +2.13%  "

$ws.Range("D33").Value = "1.422.39"
$ws.Range("E33").Value = "  -2.99%  "

$ws.Range("E34").Value = "  +2.38%  "

$ws.Range("D35").Value = "`'1.57"
$ws.Range("E35").Value = "  +2.17%  "

$ws.Range("D36").Value = "`'2.35"
$ws.Range("E36").Value = "  +0.40%  "

$ws.Range("D37").Value = "`'0.891"
$ws.Range("E37").Value = "  +2.46%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "`'0.0168"
$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "`'0.925"
$ws.Range("E39").Value = "  -1.93%  "

$ws.Range("D40").Value = "`'0.558"
$ws.Range("E40").Value = "  +1.16%  "

$ws.Range("E41").Value = "  +2.25%  "

$ws.Range("D43").Value = "`'67.00"
$ws.Range("E43").Value = "  -0.20%  "

$ws.Range("E44").Value = "  +0.43%  "

$ws.Range("E45").Value = "  +3.22%  "

$ws.Range("D46").Value = "`'1.82"
$ws.Range("E46").Value = "  +3.46%  "

$ws.Range("E47").Value = "  +0.28%  "

$ws.Range("D48").Value = "1.783.57"
$ws.Range("E48").Value = "  +1.24%  "

$ws.Range("D49").Value = "`'88.93"
$ws.Range("E49").Value = "  +2.01%  "

$ws.Range("E50").Value = "  +1.26%  "

$ws.Range("E51").Value = "  +0.68%  "
